$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rename (column headers -> snake_case field names) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Special-case full text fixes ---
$ws.Range("A332").Value = "Guanajuato"
$ws.Range("B820").Value = "Montemorelos"

# --- Title-case connector words (de/del/la/las/los/el/y) in state/municipality names ---
$ws.Range("B8").Value = "Pabellón De Arteaga"
$ws.Range("B9").Value = "Rincón De Romos"
$ws.Range("B10").Value = "San Francisco De Los Romo"
$ws.Range("B11").Value = "San José De Gracia"
$ws.Range("B31").Value = "Amatenango De La Frontera"
$ws.Range("B34").Value = "Bejucal De Ocampo"
$ws.Range("B36").Value = "Benemérito De Las Américas"
$ws.Range("B44").Value = "Comitán De Domínguez"
$ws.Range("B67").Value = "Ocozocoautla De Espinosa"
$ws.Range("B76").Value = "Salto De Agua"
$ws.Range("B77").Value = "San Cristóbal De Las Casas"
$ws.Range("B111").Value = "Guadalupe Y Calvo"
$ws.Range("B114").Value = "Hidalgo Del Parral"
$ws.Range("B127").Value = "San Francisco De Borja"
$ws.Range("B128").Value = "San Francisco Del Oro"
$ws.Range("B149").Value = "San Juan De Sabinas"
$ws.Range("B162").Value = "Villa De Álvarez"
$ws.Range("A164").Value = "Ciudad De México"
$ws.Range("B168").Value = "Cuajimalpa De Morelos"
$ws.Range("B183").Value = "Coneto De Comonfort"
$ws.Range("B197").Value = "Nombre De Dios"
$ws.Range("B201").Value = "Pánuco De Coronado"
$ws.Range("B208").Value = "San Juan Del Río"
$ws.Range("B209").Value = "San Luis Del Cordero"
$ws.Range("B210").Value = "San Pedro Del Gallo"
$ws.Range("A220").Value = "Estado De México"
$ws.Range("B220").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B223").Value = "Almoloya De Alquisiras"
$ws.Range("B224").Value = "Almoloya De Juárez"
$ws.Range("B231").Value = "Atizapán De Zaragoza"
$ws.Range("B237").Value = "Chapa De Mota"
$ws.Range("B241").Value = "Coacalco De Berriozábal"
$ws.Range("B248").Value = "Ecatepec De Morelos"
$ws.Range("B256").Value = "Ixtapan De La Sal"
$ws.Range("B257").Value = "Ixtapan Del Oro"
$ws.Range("B270").Value = "Naucalpan De Juárez"
$ws.Range("B283").Value = "San Felipe Del Progreso"
$ws.Range("B284").Value = "San Martín De Las Pirámides"
$ws.Range("B286").Value = "San Simón De Guerero"
$ws.Range("B288").Value = "Soyaniquilpan De Juárez"
$ws.Range("B297").Value = "Tenango Del Aire"
$ws.Range("B298").Value = "Tenango Del Valle"
$ws.Range("B312").Value = "Tlalnepantla De Baz"
$ws.Range("B318").Value = "Valle De Bravo"
$ws.Range("B319").Value = "Valle De Chalco Solidaridad"
$ws.Range("B320").Value = "Villa De Allende"
$ws.Range("B321").Value = "Villa Del Carbón"
$ws.Range("B335").Value = "Apaseo El Alto"
$ws.Range("B336").Value = "Apaseo El Grande"
$ws.Range("B344").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B348").Value = "Jaral Del Progreso"
$ws.Range("B356").Value = "Purísima Del Rincón"
$ws.Range("B360").Value = "San Diego De La Unión"
$ws.Range("B362").Value = "San Francisco Del Rincón"
$ws.Range("B364").Value = "San Luis De La Paz"
$ws.Range("B366").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B368").Value = "Silao De La Victoria"
$ws.Range("B373").Value = "Valle De Santiago"
$ws.Range("B379").Value = "Acapulco De Juárez"
$ws.Range("B381").Value = "Ajuchitlán Del Progreso"
$ws.Range("B382").Value = "Alcozauca De Guerero"
$ws.Range("B386").Value = "Atenango Del Río"
$ws.Range("B387").Value = "Atlamajalcingo Del Monte"
$ws.Range("B389").Value = "Atoyac De Álvarez"
$ws.Range("B390").Value = "Ayutla De Los Libres"
$ws.Range("B393").Value = "Buenavista De Cuéllar"
$ws.Range("B394").Value = "Chilapa De Álvarez"
$ws.Range("B395").Value = "Chilpancingo De Los Bravo"
$ws.Range("B396").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B401").Value = "Coyuca De Benítez"
$ws.Range("B402").Value = "Coyuca De Catalán"
$ws.Range("B406").Value = "Cuetzala Del Progreso"
$ws.Range("B407").Value = "Cutzamala De Pinzón"
$ws.Range("B413").Value = "Huitzuco De Los Figueroa"
$ws.Range("B414").Value = "Iguala De La Independencia"
$ws.Range("B416").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B417").Value = "Zihuatanejo De Azueta"
$ws.Range("B419").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B422").Value = "Mártir De Cuilapan"
$ws.Range("B435").Value = "Taxco De Alarcón"
$ws.Range("B437").Value = "Técpan De Galeana"
$ws.Range("B439").Value = "Tepecoacuilco De Trujano"
$ws.Range("B441").Value = "Tixtla De Guerero"
$ws.Range("B444").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B445").Value = "Tlapa De Comonfort"
$ws.Range("B457").Value = "Agua Blanca De Iturbide"
$ws.Range("B462").Value = "Atotonilco De Tula"
$ws.Range("B463").Value = "Atotonilco El Grande"
$ws.Range("B469").Value = "Cuautepec De Hinojosa"
$ws.Range("B474").Value = "Huasca De Ocampo"
$ws.Range("B478").Value = "Huejutla De Reyes"
$ws.Range("B481").Value = "Jacala De Ledezma"
$ws.Range("B487").Value = "Mineral Del Chico"
$ws.Range("B488").Value = "Mixquiahuala De Juárez"
$ws.Range("B490").Value = "Nopala De Villagrán"
$ws.Range("B491").Value = "Omitlán De Juárez"
$ws.Range("B492").Value = "Pachuca De Soto"
$ws.Range("B495").Value = "Progreso De Obregón"
$ws.Range("B501").Value = "Santiago De Anaya"
$ws.Range("B505").Value = "Tenango De Doria"
$ws.Range("B507").Value = "Tepehuacán De Guerero"
$ws.Range("B508").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B511").Value = "Tezontepec De Aldama"
$ws.Range("B518").Value = "Tula De Allende"
$ws.Range("B519").Value = "Tulancingo De Bravo"
$ws.Range("B522").Value = "Zacualtipán De Ángeles"
$ws.Range("B527").Value = "Ahualulco De Mercado"
$ws.Range("B532").Value = "Atemajac De Brizuela"
$ws.Range("B534").Value = "Atotonilco El Alto"
$ws.Range("B535").Value = "Autlán De Navarro"
$ws.Range("B541").Value = "Cañadas De Obregón"
$ws.Range("B547").Value = "Concepción De Buenos Aires"
$ws.Range("B548").Value = "Cuautitlán De García Barragán"
$ws.Range("B557").Value = "Encarnación De Díaz"
$ws.Range("B562").Value = "Huejuquilla El Alto"
$ws.Range("B563").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B564").Value = "Ixtlahuacán Del Río"
$ws.Range("B568").Value = "Jilotlán De Los Dolores"
$ws.Range("B574").Value = "La Manzanilla De La Paz"
$ws.Range("B575").Value = "Lagos De Moreno"
$ws.Range("B583").Value = "Ojuelos De Jalisco"
$ws.Range("B588").Value = "San Cristóbal De La Barranca"
$ws.Range("B589").Value = "San Diego De Alejandría"
$ws.Range("B591").Value = "San Juan De Los Lagos"
$ws.Range("B594").Value = "San Miguel El Alto"
$ws.Range("B595").Value = "San Sebastián Del Oeste"
$ws.Range("B596").Value = "Santa María De Los Ángeles"
$ws.Range("B597").Value = "Santa María Del Oro"
$ws.Range("B600").Value = "Talpa De Allende"
$ws.Range("B601").Value = "Tamazula De Gordiano"
$ws.Range("B604").Value = "Techaluta De Montenegro"
$ws.Range("B608").Value = "Teocuitatlán De Corona"
$ws.Range("B609").Value = "Tepatitlán De Morelos"
$ws.Range("B612").Value = "Tizapán El Alto"
$ws.Range("B613").Value = "Tlajomulco De Zúñiga"
$ws.Range("B623").Value = "Unión De San Antonio"
$ws.Range("B624").Value = "Unión De Tula"
$ws.Range("B625").Value = "Valle De Guadalupe"
$ws.Range("B626").Value = "Valle De Juárez"
$ws.Range("B631").Value = "Yahualica De González Gallo"
$ws.Range("B632").Value = "Zacoalco De Torres"
$ws.Range("B635").Value = "Zapotlán Del Rey"
$ws.Range("B636").Value = "Zapotlán El Grande"
$ws.Range("B662").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B664").Value = "Cojumatlán De Régules"
$ws.Range("B729").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B755").Value = "Coatlán Del Río"
$ws.Range("B763").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B767").Value = "Puente De Ixtla"
$ws.Range("B773").Value = "Tetela Del Volcán"
$ws.Range("B774").Value = "Tlaltizapán De Zapata"
$ws.Range("B782").Value = "Zacualpan De Amilpas"
$ws.Range("B785").Value = "Bahía De Banderas"
$ws.Range("B788").Value = "Ixtlán Del Río"
$ws.Range("B795").Value = "Santa María Del Oro"
$ws.Range("B808").Value = "Ciénega De Flores"
$ws.Range("B823").Value = "San Nicolás De Los Garza"
$ws.Range("B827").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B833").Value = "Chalcatongo De Hidalgo"
$ws.Range("B834").Value = "Ciénega De Zimatlán"
$ws.Range("B839").Value = "Cuilápam De Guerero"
$ws.Range("B840").Value = "Cuyamecalco Villa De Zaragoza"
$ws.Range("B841").Value = "El Barrio De La Soledad"
$ws.Range("B843").Value = "Guevea De Humboldt"
$ws.Range("B844").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B845").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B846").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B847").Value = "Huautla De Jiménez"
$ws.Range("B848").Value = "Ixtlán De Juárez"
$ws.Range("B849").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B857").Value = "Mariscala De Juárez"
$ws.Range("B858").Value = "Mártires De Tacubaya"
$ws.Range("B860").Value = "Mazatlán Villa De Flores"
$ws.Range("B862").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B864").Value = "Oaxaca De Juárez"
$ws.Range("B865").Value = "Ocotlán De Morelos"
$ws.Range("B866").Value = "Pinotepa De Don Luis"
$ws.Range("B868").Value = "Putla Villa De Guerero"
$ws.Range("B869").Value = "Reforma De Pineda"
$ws.Range("B883").Value = "San Antonio De La Cal"
$ws.Range("B885").Value = "San Baltazar Yatzachi El Bajo"
$ws.Range("B890").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B908").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B917").Value = "San Juan De Los Cués"
$ws.Range("B943").Value = "San Mateo Del Mar"
$ws.Range("B951").Value = "San Miguel Del Puerto"
$ws.Range("B952").Value = "San Miguel El Grande"
$ws.Range("B964").Value = "San Pablo Villa De Mitla"
$ws.Range("B969").Value = "San Pedro El Alto"
$ws.Range("B983").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B984").Value = "San Pedro Y San Pablo Tequixtepec"
$ws.Range("B993").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B1006").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B1042").Value = "Santo Domingo De Morelos"
$ws.Range("B1056").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B1057").Value = "Tataltepec De Valdés"
$ws.Range("B1058").Value = "Teotitlán De Flores Magón"
$ws.Range("B1060").Value = "Tepelmeme Villa De Morelos"
$ws.Range("B1061").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B1062").Value = "Tlacolula De Matamoros"
$ws.Range("B1064").Value = "Villa De Etla"
$ws.Range("B1065").Value = "Villa De Tamazulápam Del Progreso"
$ws.Range("B1066").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B1067").Value = "Villa Sola De Vega"
$ws.Range("B1069").Value = "Zapotitlán Del Río"
$ws.Range("B1072").Value = "Zimatlán De Álvarez"
$ws.Range("B1091").Value = "Ayotoxco De Guerero"
$ws.Range("B1096").Value = "Chalchicomula De Sesma"
$ws.Range("B1102").Value = "Chila De La Sal"
$ws.Range("B1112").Value = "Cuapiaxtla De Madero"
$ws.Range("B1115").Value = "Cuayuca De Andrade"
$ws.Range("B1116").Value = "Cuetzalan Del Progreso"
$ws.Range("B1131").Value = "Huehuetlán El Chico"
$ws.Range("B1132").Value = "Huehuetlán El Grande"
$ws.Range("B1136").Value = "Huitzilan De Serdán"
$ws.Range("B1137").Value = "Ixcamilpa De Guerero"
$ws.Range("B1140").Value = "Izúcar De Matamoros"
$ws.Range("B1150").Value = "Los Reyes De Juárez"
$ws.Range("B1151").Value = "Mazapiltepec De Juárez"
$ws.Range("B1160").Value = "Palmar De Bravo"
$ws.Range("B1180").Value = "San Nicolás De Los Ranchos"
$ws.Range("B1184").Value = "San Salvador El Seco"
$ws.Range("B1185").Value = "San Salvador El Verde"
$ws.Range("B1192").Value = "Tecali De Herrera"
$ws.Range("B1200").Value = "Tepanco De López"
$ws.Range("B1201").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B1206").Value = "Tepexi De Rodríguez"
$ws.Range("B1208").Value = "Tetela De Ocampo"
$ws.Range("B1213").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B1244").Value = "Amealco De Bonfil"
$ws.Range("B1245").Value = "Cadereyta De Montes"
$ws.Range("B1248").Value = "Jalpan De Serra"
$ws.Range("B1249").Value = "Landa De Matamoros"
$ws.Range("B1251").Value = "Pinal De Amoles"
$ws.Range("B1254").Value = "San Juan Del Río"
$ws.Range("B1265").Value = "Armadillo De Los Infante"
$ws.Range("B1271").Value = "Ciudad Del Maíz"
$ws.Range("B1279").Value = "Mexquitic De Carmona"
$ws.Range("B1285").Value = "San Ciro De Acosta"
$ws.Range("B1288").Value = "Santa María Del Río"
$ws.Range("B1290").Value = "Soledad De Graciano Sánchez"
$ws.Range("B1297").Value = "Tanquián De Escobedo"
$ws.Range("B1301").Value = "Villa De Arista"
$ws.Range("B1302").Value = "Villa De Arriaga"
$ws.Range("B1303").Value = "Villa De Ramos"
$ws.Range("B1304").Value = "Villa De Reyes"
$ws.Range("B1335").Value = "Nacozari De García"
$ws.Range("B1340").Value = "San Felipe De Jesús"
$ws.Range("B1379").Value = "Soto La Marina"
$ws.Range("B1386").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B1388").Value = "Amaxac De Guerero"
$ws.Range("B1389").Value = "Apetatitlán De Antonio Carvajal"
$ws.Range("B1395").Value = "Contla De Juan Cuamatzi"
$ws.Range("B1400").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B1403").Value = "Mazatecochco De José María Morelos"
$ws.Range("B1404").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B1409").Value = "San Pablo Del Monte"
$ws.Range("B1414").Value = "Tepetitla De Lardizábal"
$ws.Range("B1417").Value = "Tetla De La Solidaridad"
$ws.Range("B1437").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B1440").Value = "Amatlán De Los Reyes"
$ws.Range("B1449").Value = "Boca Del Río"
$ws.Range("B1452").Value = "Castillo De Teayo"
$ws.Range("B1454").Value = "Cazones De Herrera"
$ws.Range("B1470").Value = "Cosamaloapan De Carpio"
$ws.Range("B1471").Value = "Cosautlán De Carvajal"
$ws.Range("B1486").Value = "Hueyapan De Ocampo"
$ws.Range("B1487").Value = "Ignacio De La Llave"
$ws.Range("B1490").Value = "Ixhuacán De Los Reyes"
$ws.Range("B1491").Value = "Ixhuatlán De Madero"
$ws.Range("B1492").Value = "Ixhuatlán Del Café"
$ws.Range("B1493").Value = "Ixhuatlán Del Sureste"
$ws.Range("B1502").Value = "Juchique De Ferrer"
$ws.Range("B1505").Value = "Landero Y Coss"
$ws.Range("B1507").Value = "Las Vigas De Ramírez"
$ws.Range("B1508").Value = "Lerdo De Tejada"
$ws.Range("B1511").Value = "Martínez De La Torre"
$ws.Range("B1517").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B1526").Value = "Ozuluama De Mascareñas"
$ws.Range("B1529").Value = "Paso De Ovejas"
$ws.Range("B1530").Value = "Paso Del Macho"
$ws.Range("B1534").Value = "Poza Rica De Hidalgo"
$ws.Range("B1543").Value = "Sayula De Alemán"
$ws.Range("B1545").Value = "Soledad De Doblado"
$ws.Range("B1566").Value = "Tlacotepec De Mejía"
$ws.Range("B1575").Value = "Vega De Alatorre"
$ws.Range("B1586").Value = "Zontecomatlán De López Y Fuentes"
$ws.Range("B1587").Value = "Zozocolco De Hidalgo"
$ws.Range("B1601").Value = "Cañitas De Felipe Pescador"
$ws.Range("B1603").Value = "Concepción Del Oro"
$ws.Range("B1605").Value = "El Plateado De Joaquín Amaro"
$ws.Range("B1615").Value = "Jiménez Del Teul"
$ws.Range("B1621").Value = "Mezquital Del Oro"
$ws.Range("B1626").Value = "Nochistlán De Mejía"
$ws.Range("B1627").Value = "Noria De Ángeles"
$ws.Range("B1638").Value = "Teúl De González Ortega"
$ws.Range("B1639").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1640").Value = "Trinidad García De La Cadena"
$ws.Range("B1643").Value = "Villa De Cos"

# --- Minor floating point literal refresh on percentage column ---
$ws.Range("D176").Value = 0.0009750297265160524
$ws.Range("D551").Value = 0.0009750297265160524
$ws.Range("D644").Value = 0.0009750297265160524
$ws.Range("D696").Value = 0.0009750297265160524
$ws.Range("D798").Value = 0.0009036860879904876
$ws.Range("D1039").Value = 0.0009036860879904876
$ws.Range("D1273").Value = 0.0009036860879904876
$ws.Range("D1315").Value = 0.0009036860879904876
$ws.Range("D1380").Value = 0.0009036860879904876
$ws.Range("D1468").Value = 0.0009036860879904876
$ws.Range("D1524").Value = 0.0009036860879904876
$ws.Range("D1633").Value = 0.0009036860879904876

# --- Delete trailing footer/metadata rows (1652-1656) ---
$ws.Range("A1652:A1656").EntireRow.Delete() | Out-Null

Write-Host "Edit complete"
